{"js": "// \"fixed all list and new tag associations\"\n//\n// Repairs the List1 / List6 / List7 / List8 / List3change / List4change\n// paragraph (and linked character) styles:\n//   - List1, List6, List8, List3change, List4change: turn OFF contextual\n//     spacing (List7 keeps its contextual spacing untouched).\n//   - List6, List7, List8: detach from the \"Heading 3\" base style and make\n//     their own Times New Roman / 12pt (+ complex-script) font explicit\n//     instead of inheriting it, and drop the outline level (these are\n//     plain list styles, not TOC-eligible headings).\n//   - List6, List8: no longer force bold off via an inherited heading\n//     (now simply not-bold on their own properties).\n//   - List7 / List7Char: become bold.\n\nconst styles = context.document.getStyles();\n\nconst list1 = styles.getByNameOrNullObject(\"List 1\");\nconst list6 = styles.getByNameOrNullObject(\"List 6\");\nconst list7 = styles.getByNameOrNullObject(\"List 7\");\nconst list7Char = styles.getByNameOrNullObject(\"List 7 Char\");\nconst list8 = styles.getByNameOrNullObject(\"List 8\");\nconst list3change = styles.getByNameOrNullObject(\"List 3_change\");\nconst list4change = styles.getByNameOrNullObject(\"List 4_change\");\n\nawait context.sync();\n\n// --- List 1 ------------------------------------------------------------\nif (!list1.isNullObject) {\n  list1.noSpaceBetweenParagraphsOfSameStyle = false;\n}\n\n// --- List 6 --------------------------------------------------------------\nif (!list6.isNullObject) {\n  list6.baseStyle = \"\";\n  list6.noSpaceBetweenParagraphsOfSameStyle = false;\n  list6.paragraphFormat.outlineLevel = 10; // wdOutlineLevelBodyText (no heading outline level)\n  list6.font.name = \"Times New Roman\";\n  list6.font.nameAscii = \"Times New Roman\";\n  list6.font.nameBidirectional = \"Times New Roman\";\n  list6.font.bold = false;\n  list6.font.size = 12;\n  list6.font.sizeBidirectional = 12;\n}\n\n// --- List 7 (and its linked character style) ----------------------------\nif (!list7.isNullObject) {\n  list7.baseStyle = \"\";\n  list7.paragraphFormat.outlineLevel = 10;\n  list7.font.name = \"Times New Roman\";\n  list7.font.nameAscii = \"Times New Roman\";\n  list7.font.nameBidirectional = \"Times New Roman\";\n  list7.font.bold = true;\n  list7.font.size = 12;\n  list7.font.sizeBidirectional = 12;\n}\nif (!list7Char.isNullObject) {\n  list7Char.font.bold = true;\n}\n\n// --- List 8 ----------------------------------------------------------------\nif (!list8.isNullObject) {\n  list8.baseStyle = \"\";\n  list8.noSpaceBetweenParagraphsOfSameStyle = false;\n  list8.paragraphFormat.outlineLevel = 10;\n  list8.font.name = \"Times New Roman\";\n  list8.font.nameAscii = \"Times New Roman\";\n  list8.font.nameBidirectional = \"Times New Roman\";\n  list8.font.bold = false;\n  list8.font.size = 12;\n  list8.font.sizeBidirectional = 12;\n}\n\n// --- List 3_change / List 4_change ------------------------------------------\nif (!list3change.isNullObject) {\n  list3change.noSpaceBetweenParagraphsOfSameStyle = false;\n}\nif (!list4change.isNullObject) {\n  list4change.noSpaceBetweenParagraphsOfSameStyle = false;\n}\n\nawait context.sync();\n", "ps1": "# \"fixed all list and new tag associations\"\n#\n# Repairs the List1/List6/List7/List8/List3change/List4change paragraph\n# (and linked character) styles:\n#   - List1, List6, List8, List3change, List4change: turn OFF contextual\n#     spacing (List7 keeps it).\n#   - List6, List7, List8: detach from the \"Heading 3\" base style and\n#     make their own Times New Roman / 12pt (+ complex-script) font\n#     explicit instead of inheriting it, and drop the outline level\n#     (these are plain list styles, not TOC-eligible headings).\n#   - List6, List8: no longer force bold off via an inherited heading\n#     (now simply not-bold on their own rPr).\n#   - List7 / List7Char: become bold.\n\n$d = $word.ActiveDocument\n\n# --- List 1 --------------------------------------------------------------\n$list1 = $d.Styles(\"List1\")\n$list1.NoSpaceBetweenParagraphsOfSameStyle = $false\n\n# --- List 6 ----------------------------------------------------------------\n$list6 = $d.Styles(\"List6\")\n$list6.BaseStyle = \"\"\n$list6.NoSpaceBetweenParagraphsOfSameStyle = $false\n$list6.ParagraphFormat.OutlineLevel = 10\n$list6.Font.Name = \"Times New Roman\"\n$list6.Font.NameAscii = \"Times New Roman\"\n$list6.Font.NameBi = \"Times New Roman\"\n$list6.Font.Bold = $false\n$list6.Font.Size = 12\n$list6.Font.SizeBi = 12\n\n# --- List 7 (and its linked character style) --------------------------------\n$list7 = $d.Styles(\"List7\")\n$list7.BaseStyle = \"\"\n$list7.ParagraphFormat.OutlineLevel = 10\n$list7.Font.Name = \"Times New Roman\"\n$list7.Font.NameAscii = \"Times New Roman\"\n$list7.Font.NameBi = \"Times New Roman\"\n$list7.Font.Bold = $true\n$list7.Font.Size = 12\n$list7.Font.SizeBi = 12\n\n$list7Char = $d.Styles(\"List7Char\")\n$list7Char.Font.Bold = $true\n\n# --- List 8 ------------------------------------------------------------------\n$list8 = $d.Styles(\"List8\")\n$list8.BaseStyle = \"\"\n$list8.NoSpaceBetweenParagraphsOfSameStyle = $false\n$list8.ParagraphFormat.OutlineLevel = 10\n$list8.Font.Name = \"Times New Roman\"\n$list8.Font.NameAscii = \"Times New Roman\"\n$list8.Font.NameBi = \"Times New Roman\"\n$list8.Font.Bold = $false\n$list8.Font.Size = 12\n$list8.Font.SizeBi = 12\n\n# --- List 3_change / List 4_change -------------------------------------------\n$list3change = $d.Styles(\"List3change\")\n$list3change.NoSpaceBetweenParagraphsOfSameStyle = $false\n\n$list4change = $d.Styles(\"List4change\")\n$list4change.NoSpaceBetweenParagraphsOfSameStyle = $false\n"}
